$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the "Valor Mora" amounts between period 2506 (row 16) and period 2412 (row 22)
$ws.Range("F16").Value = 138000
$ws.Range("F22").Value = 119600
